$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-07 Saturday", 2)
$d.Content.Find.Execute("17÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷9=", 2)
$d.Content.Find.Execute("92÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷6=", 2)
$d.Content.Find.Execute("66÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=", 2)
$d.Content.Find.Execute("95÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷2=", 2)
$d.Content.Find.Execute("55÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷6=", 2)
$d.Content.Find.Execute("91÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=", 2)
$d.Content.Find.Execute("54÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷3=", 2)
$d.Content.Find.Execute("49÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=", 2)
$d.Content.Find.Execute("79÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=", 2)
$d.Content.Find.Execute("84÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 2)
$d.Content.Find.Execute("29÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷3=", 2)
$d.Content.Find.Execute("37÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=", 2)
$d.Content.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=", 2)
$d.Content.Find.Execute("11÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷7=", 2)
$d.Content.Find.Execute("48÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=", 2)
$d.Content.Find.Execute("55÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷8=", 2)
$d.Content.Find.Execute("48÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷5=", 2)
$d.Content.Find.Execute("36÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=", 2)
$d.Content.Find.Execute("62÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=", 2)
$d.Content.Find.Execute("15÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=", 2)
$d.Content.Find.Execute("71÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷9=", 2)
$d.Content.Find.Execute("90÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=", 2)
$d.Content.Find.Execute("97÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷2=", 2)
$d.Content.Find.Execute("80÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=", 2)
$d.Content.Find.Execute("20÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=", 2)
